$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value, forcing the correct data type (numeric vs text),
# and if the cell's underlying style needs to change (text <-> numeric format),
# copy the number format from a same-format anchor cell elsewhere on the sheet
# (format-only paste leaves the just-assigned value untouched).
function Set-Cell($addr, $val, $kind, $styleAnchor) {
    $c = $ws.Range($addr)
    if ($kind -eq "s") {
        # Force text interpretation (so numeric-looking strings like "0" are
        # stored as shared-string text, not converted to a number).
        $c.NumberFormat = "@"
        $c.Value = $val
    } else {
        $c.Value = $val
    }
    if ($styleAnchor) {
        $src = $ws.Range($styleAnchor)
        $src.Copy()
        $c.PasteSpecial(-4122)   # xlPasteFormats: formats only, value is preserved
    }
}

# --- Header: Volume/issue number and report date range ---
# "Volume 31   Number  20" -> "Volume 31   Number  21"
$ws.Range("A8").Characters(21, 2).Text = "21"

# "Report Covering the Week  5/13/2024  Through  5/19/2024"
#   -> "Report Covering the Week  5/20/2024  Through  5/26/2024"
$ws.Range("C9").Characters(27, 9).Text = "5/20/2024"
$ws.Range("C9").Characters(47, 9).Text = "5/26/2024"

# --- Weekly crime statistics table (rows 15-30) ---
    Set-Cell "C15" 1 "n" "F14"
    Set-Cell "D15" "0" "s" "A14"
    Set-Cell "E15" "***.*" "s" "A14"
    Set-Cell "F15" 1 "n" "F14"
    Set-Cell "H15" -75 "n" $null
    Set-Cell "I15" 14 "n" $null
    Set-Cell "K15" 40 "n" $null
    Set-Cell "L15" 133.333333333333 "n" $null
    Set-Cell "M15" 40 "n" $null
    Set-Cell "N15" 16.666666666666 "n" $null
    Set-Cell "C16" 4 "n" $null
    Set-Cell "D16" 16 "n" $null
    Set-Cell "E16" -75 "n" $null
    Set-Cell "F16" 12 "n" $null
    Set-Cell "G16" 28 "n" $null
    Set-Cell "H16" -57.142857142857 "n" $null
    Set-Cell "I16" 106 "n" $null
    Set-Cell "J16" 101 "n" $null
    Set-Cell "K16" 4.950495049504 "n" $null
    Set-Cell "L16" 20.454545454545 "n" $null
    Set-Cell "M16" -1.851851851851 "n" $null
    Set-Cell "N16" -64.548494983277 "n" $null
    Set-Cell "C17" 11 "n" $null
    Set-Cell "D17" 7 "n" $null
    Set-Cell "E17" 57.142857142857 "n" $null
    Set-Cell "F17" 41 "n" $null
    Set-Cell "G17" 38 "n" $null
    Set-Cell "H17" 7.894736842105 "n" $null
    Set-Cell "I17" 176 "n" $null
    Set-Cell "J17" 161 "n" $null
    Set-Cell "K17" 9.316770186335 "n" $null
    Set-Cell "L17" 32.330827067669 "n" $null
    Set-Cell "M17" 79.591836734693 "n" $null
    Set-Cell "N17" 60 "n" $null
    Set-Cell "C18" 6 "n" $null
    Set-Cell "E18" 0 "n" $null
    Set-Cell "F18" 17 "n" $null
    Set-Cell "H18" 21.428571428571 "n" $null
    Set-Cell "I18" 75 "n" $null
    Set-Cell "J18" 95 "n" $null
    Set-Cell "K18" -21.052631578947 "n" $null
    Set-Cell "L18" 41.509433962264 "n" $null
    Set-Cell "M18" -43.609022556391 "n" $null
    Set-Cell "N18" -88.114104595879 "n" $null
    Set-Cell "C19" 8 "n" $null
    Set-Cell "D19" 12 "n" $null
    Set-Cell "E19" -33.333333333333 "n" $null
    Set-Cell "F19" 58 "n" $null
    Set-Cell "G19" 52 "n" $null
    Set-Cell "H19" 11.538461538461 "n" $null
    Set-Cell "I19" 364 "n" $null
    Set-Cell "J19" 245 "n" $null
    Set-Cell "K19" 48.571428571428 "n" $null
    Set-Cell "L19" 38.403041825095 "n" $null
    Set-Cell "M19" 142.666666666667 "n" $null
    Set-Cell "N19" 62.5 "n" $null
    Set-Cell "C20" 7 "n" $null
    Set-Cell "D20" 19 "n" $null
    Set-Cell "E20" -63.157894736842 "n" $null
    Set-Cell "F20" 33 "n" $null
    Set-Cell "G20" 64 "n" $null
    Set-Cell "H20" -48.4375 "n" $null
    Set-Cell "I20" 176 "n" $null
    Set-Cell "J20" 208 "n" $null
    Set-Cell "K20" -15.384615384615 "n" $null
    Set-Cell "L20" 25.714285714285 "n" $null
    Set-Cell "M20" 91.304347826087 "n" $null
    Set-Cell "N20" -76.564580559254 "n" $null
    Set-Cell "C21" 37 "n" $null
    Set-Cell "D21" 60 "n" $null
    Set-Cell "E21" -38.333333333333 "n" $null
    Set-Cell "F21" 163 "n" $null
    Set-Cell "G21" 200 "n" $null
    Set-Cell "H21" -18.5 "n" $null
    Set-Cell "I21" 913 "n" $null
    Set-Cell "J21" 822 "n" $null
    Set-Cell "K21" 11.070559610705 "n" $null
    Set-Cell "L21" 33.090379008746 "n" $null
    Set-Cell "M21" 53.703703703703 "n" $null
    Set-Cell "N21" -55.068897637795 "n" $null
    Set-Cell "D22" 1 "n" "F14"
    Set-Cell "E22" -100 "n" "K14"
    Set-Cell "G22" 1 "n" "F14"
    Set-Cell "H22" -100 "n" "K14"
    Set-Cell "J22" 8 "n" $null
    Set-Cell "K22" 12.5 "n" $null
    Set-Cell "C23" 5 "n" $null
    Set-Cell "D23" 1 "n" $null
    Set-Cell "E23" 400 "n" $null
    Set-Cell "F23" 10 "n" $null
    Set-Cell "G23" 9 "n" $null
    Set-Cell "H23" 11.111111111111 "n" $null
    Set-Cell "I23" 49 "n" $null
    Set-Cell "J23" 54 "n" $null
    Set-Cell "K23" -9.259259259259 "n" $null
    Set-Cell "L23" 19.512195121951 "n" $null
    Set-Cell "M23" 113.04347826087 "n" $null
    Set-Cell "C24" 30 "n" $null
    Set-Cell "D24" 40 "n" $null
    Set-Cell "E24" -25 "n" $null
    Set-Cell "F24" 113 "n" $null
    Set-Cell "G24" 120 "n" $null
    Set-Cell "H24" -5.833333333333 "n" $null
    Set-Cell "I24" 593 "n" $null
    Set-Cell "J24" 615 "n" $null
    Set-Cell "K24" -3.577235772357 "n" $null
    Set-Cell "L24" 20.773930753564 "n" $null
    Set-Cell "M24" 69.428571428571 "n" $null
    Set-Cell "C25" 12 "n" $null
    Set-Cell "D25" 19 "n" $null
    Set-Cell "E25" -36.842105263157 "n" $null
    Set-Cell "F25" 34 "n" $null
    Set-Cell "G25" 44 "n" $null
    Set-Cell "H25" -22.727272727272 "n" $null
    Set-Cell "I25" 246 "n" $null
    Set-Cell "J25" 252 "n" $null
    Set-Cell "K25" -2.380952380952 "n" $null
    Set-Cell "L25" 30.851063829787 "n" $null
    Set-Cell "C26" 14 "n" $null
    Set-Cell "D26" 7 "n" $null
    Set-Cell "E26" 100 "n" $null
    Set-Cell "F26" 53 "n" $null
    Set-Cell "G26" 46 "n" $null
    Set-Cell "H26" 15.217391304347 "n" $null
    Set-Cell "I26" 214 "n" $null
    Set-Cell "J26" 226 "n" $null
    Set-Cell "K26" -5.309734513274 "n" $null
    Set-Cell "L26" 7 "n" $null
    Set-Cell "M26" -12.653061224489 "n" $null
    Set-Cell "C27" 2 "n" $null
    Set-Cell "D27" "0" "s" "A14"
    Set-Cell "E27" "***.*" "s" "A14"
    Set-Cell "I27" 17 "n" $null
    Set-Cell "K27" 13.333333333333 "n" $null
    Set-Cell "L27" 13.333333333333 "n" $null
    Set-Cell "D28" 2 "n" $null
    Set-Cell "E28" -50 "n" $null
    Set-Cell "G28" 9 "n" $null
    Set-Cell "H28" -55.555555555555 "n" $null
    Set-Cell "I28" 25 "n" $null
    Set-Cell "J28" 25 "n" $null
    Set-Cell "K28" 0 "n" $null
    Set-Cell "L28" 31.578947368421 "n" $null
    Set-Cell "F29" 1 "n" $null
    Set-Cell "H29" -50 "n" $null
    Set-Cell "N29" -72.222222222222 "n" $null
    Set-Cell "F30" 1 "n" $null
    Set-Cell "H30" 0 "n" $null
    Set-Cell "N30" -66.666666666666 "n" $null
